# Auto-generated edit script applying the diff (F-column count bumps on
# sheets 展览/演出/全部类型, plus a content shift + new row on 全部类型 rows 25-32).
$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 145
$ws.Range("F3").Value = 419
$ws.Range("F5").Value = 47
$ws.Range("F6").Value = 1274
$ws.Range("F7").Value = 469
$ws.Range("F8").Value = 102
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 192
$ws.Range("F12").Value = 1071
$ws.Range("F14").Value = 277
$ws.Range("F15").Value = 217
$ws.Range("F16").Value = 1573
$ws.Range("F17").Value = 577
$ws.Range("F19").Value = 371
$ws.Range("F21").Value = 886
$ws.Range("F22").Value = 1176
$ws.Range("F24").Value = 1910
$ws.Range("F25").Value = 2716
$ws.Range("F26").Value = 1499
$ws.Range("F28").Value = 70
$ws.Range("F29").Value = 505
$ws.Range("F30").Value = 827
$ws.Range("F31").Value = 1433
$ws.Range("F33").Value = 1498
$ws.Range("F34").Value = 175
$ws.Range("F36").Value = 806
$ws.Range("F37").Value = 698
$ws.Range("F38").Value = 715
$ws.Range("F39").Value = 913
$ws.Range("F40").Value = 381
$ws.Range("F41").Value = 273

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 730
$ws.Range("F23").Value = 26

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 145
$ws.Range("F5").Value = 419
$ws.Range("F7").Value = 47
$ws.Range("F10").Value = 1274
$ws.Range("F11").Value = 469
$ws.Range("F12").Value = 102
$ws.Range("F14").Value = 159
$ws.Range("F15").Value = 192
$ws.Range("F16").Value = 1071
$ws.Range("F18").Value = 277
$ws.Range("F20").Value = 217
$ws.Range("F21").Value = 1573
$ws.Range("F22").Value = 577
$ws.Range("F24").Value = 371
$ws.Range("B25").Value = "'2024-06-07"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "杭州·【鼓楼西戏剧】周一围领衔主演·《枕头人》10周年纪念版"
$ws.Range("D25").Value = "杭州市江干区新业路39号 杭州大剧院"
$ws.Range("E25").Value = "2024.06.07 19:30-06.08 22:00"
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 480
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=84902"
$ws.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202404/VZQS8SJP1714020772683.jpeg"
$ws.Range("B26").Value = "'2024-06-08"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "杭州·第38届漫展x原崩铁only"
$ws.Range("D26").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E26").Value = "2024.06.08 10:30-06.09 17:00"
$ws.Range("F26").Value = 1176
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=84802"
$ws.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202404/QsVzW6XP1712908414935.jpeg"
$ws.Range("B27").Value = "'2024-06-09"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "杭州·第三届日夜国乙only"
$ws.Range("D27").Value = "创意路1号 中国智谷富春园区"
$ws.Range("E27").Value = "2024.06.09 10:00-06.09 23:00"
$ws.Range("F27").Value = 2716
$ws.Range("G27").Value = 58
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=82618"
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png"
$ws.Range("B28").Value = "'2024-06-14"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "杭州·苗阜王声 青曲社相声全国巡演"
$ws.Range("D28").Value = "湖墅南路138号 杭州浙话艺术剧院"
$ws.Range("E28").Value = "2024.06.14 19:30-06.14 22:00"
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 280
$ws.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=83382"
$ws.Range("I28").Value = "//i1.hdslb.com/bfs/openplatform/202403/hUGL3xz01711346789039.jpeg"
$ws.Range("B29").Value = "'2024-06-15"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "杭州·次元盛典1.0"
$ws.Range("D29").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E29").Value = "2024.06.15 10:00-06.16 17:00"
$ws.Range("F29").Value = 1499
$ws.Range("G29").Value = 68
$ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=83672"
$ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202404/yZAi07mM1712033477653.jpeg"
$ws.Range("C30").Value = "杭州·第三届动漫迷城嘉年华·毕业泳池"
$ws.Range("D30").Value = "东新路21号 九龙仓君玺"
$ws.Range("E30").Value = "2024.06.15 10:00-06.15 17:00"
$ws.Range("F30").Value = 72
$ws.Range("G30").Value = 70
$ws.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=84338"
$ws.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202404/wQAlXTnK1713202337669.jpeg"
$ws.Range("B31").Value = "'2024-06-22"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "杭州·巅峰对决·排球少年ONLY"
$ws.Range("D31").Value = "金桥北路990号 万达广场(杭州富阳店)"
$ws.Range("E31").Value = "2024.06.22 10:00-06.22 17:00"
$ws.Range("F31").Value = 70
$ws.Range("G31").Value = 60
$ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=85095"
$ws.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202404/3WtpRjjo1714197500930.jpeg"
$ws.Range("C32").Value = "杭州·时光代理人「惊喜节拍」LIVE｜2024音乐巡演"
$ws.Range("D32").Value = "新北街85号三层G2-302 杭州大麦66 LIVEHOUSE"
$ws.Range("E32").Value = "2024.06.22 20:00-06.22 22:00"
$ws.Range("F32").Value = 730
$ws.Range("G32").Value = 399
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=85043"
$ws.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202404/3nufasTp1714404961103.jpeg"
$ws.Range("F34").Value = 505
$ws.Range("F35").Value = 827
$ws.Range("F36").Value = 1433
$ws.Range("F40").Value = 1498
$ws.Range("F41").Value = 806
$ws.Range("F42").Value = 698
$ws.Range("F43").Value = 715
$ws.Range("F44").Value = 913
$ws.Range("F45").Value = 381
$ws.Range("F46").Value = 26
$ws.Range("F48").Value = 273

